$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 191
$ws.Range("I2").Value = 551
$ws.Range("J2").Value = 2201
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 618
$ws.Range("M2").Value = 37
$ws.Range("N2").Value = 413
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 229
$ws.Range("T2").Value = 392
$ws.Range("U2").Value = 37
$ws.Range("V2").Value = 3420
$ws.Range("X2").Value = 3346
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 15
